# Add the "funding graph" backend columns (H:J) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (H1:J1) - shared strings are appended in this exact order.
$ws.Range("H1").Value = "Total Funding"
$ws.Range("I1").Value = "Monthly Budget"
$ws.Range("J1").Value = "Expenditure"

# New description cell under the "Total Funding" header.
$ws.Range("H2").Value = "Each cell represents budget for that month"

# Match the bold header formatting already used by A1:G1.
$ws.Range("H1:J1").Font.Bold = $true

# Resize the columns (G grew because its header text wraps, H:J are new).
$ws.Columns.Item(7).ColumnWidth = 31.333333333333332
$ws.Columns.Item(8).ColumnWidth = 11.666666666666666
$ws.Columns.Item(9).ColumnWidth = 14
$ws.Columns.Item(10).ColumnWidth = 10.333333333333334

# Move the active selection to K2, matching the saved view state.
$ws.Range("K2").Select() | Out-Null
